$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the title (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Fruit Shop for Free - Review of
#    Gameplay and Features" right before the final paragraph (the one that
#    currently holds the italic image-prompt text).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Shop for Free - Review of Gameplay and Features</w:t></w:r></w:p>' +
            '<w:p/>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($newParaXml) | Out-Null

# InsertXML above leaves behind a stray empty paragraph between the newly
# inserted bold paragraph and the original (untouched) final paragraph;
# remove it so the final paragraph is unchanged and directly follows.
$strayIndex = $d.Paragraphs.Count - 1
$strayPara = $d.Paragraphs.Item($strayIndex)
if ($strayPara.Range.Text -eq [string][char]13) {
    $strayPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Swap the old italic image-generation prompt for the new meta
#    description text, keeping the paragraph's italic formatting intact.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a cartoon-style image featuring a happy Maya warrior with glasses for the game " + [char]34 + "Fruit Shop" + [char]34 + ". The image should showcase the Maya warrior holding a basket of colorful fruits, with a backdrop of a fruit stand. The Maya warrior should have a joyful expression on his face, indicative of a successful day at the fruit stand. The fruits in the basket should be easily recognizable and brightly colored, including apples, oranges, cherries, grapes, and watermelon. The image should be vibrant and eye-catching, making it relevant to gamers of all ages.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Fruit Shop and play it for free. Learn about its gameplay, pay lines, free spins, wild feature and Return to Player value.",
    2) | Out-Null
